$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = "Semester 1 2016"

$ws.Range("E11").Value = 4

$ws.Range("L14").Value = $null
$ws.Range("M14").Value = 10

$ws.Range("F16").Value = 10
$ws.Range("N16").Value = 10
$ws.Range("O16").Value = 10

$ws.Range("J18").Value = $null
$ws.Range("K18").Value = 10
$ws.Range("L18").Value = 10

$ws.Range("F20").Value = $null
$ws.Range("J20").Value = 10

$ws.Range("P22").Value = 10
$ws.Range("Q22").Value = 10

$ws.Range("Q23").Select()
